$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9 (added first so new shared strings land in the order the
#     target workbook uses: Jeremy Faludi, HTML, 6 mo-1 year overall,
#     Once every 10 years are the first new unique strings) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Jeremy Faludi"
$ws.Range("C9").Value = "HTML"
$ws.Range("D9").Value = "x"
$ws.Range("E9").Value = "x"
$ws.Range("G9").Value = "x"
$ws.Range("M9").Value = "6 mo-1 year overall"
$ws.Range("Q9").Value = "Once every 10 years"

# --- Row 7 ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Daniel Seita"
$ws.Range("C7").Value = "Python"

# --- Row 8 ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Christine Gregg"
$ws.Range("C8").Value = "MATLAB"

# Update the active selection to match the edited workbook (was S7, now C8)
$ws.Range("C8").Select()
